# Data_Extract_From_World_Development_Indicators.xlsx
# Country list changes: "Algeria" -> "Colombia" and "South Africa" -> "Kenya",
# each with refreshed population figures on the "Data" sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("Series - Metadata")

# Row 12: Algeria -> Colombia
$ws1.Range("A12").Value = "Colombia"
$ws1.Range("B12").Value = 34539596
$ws1.Range("C12").Value = 4412670

# Row 15: South Africa -> Kenya
$ws1.Range("A15").Value = "Kenya"
$ws1.Range("B15").Value = 30694157
$ws1.Range("C15").Value = 1274302

# Update the saved selections for each sheet (set the non-active sheet first
# so the final Select() on the Data sheet leaves it as the active tab).
$ws2.Range("A5").Select() | Out-Null
$ws1.Range("D21").Select() | Out-Null
